# FInd Empty Coulmn Automatically
#
# Scan the Day1..Day40 attendance columns (E..J used here) for rows 2-6
# and automatically find empty columns/cells, filling them in with the
# appropriate Present/Absent attendance value. A couple of previously
# mis-marked Day1 cells are also refreshed with their corrected value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Attendance values to apply to specific (row, column) cells that are
# currently empty. Cells not listed here, and cells that already contain
# a value, are left untouched.
$attendance = @{
    2 = @{ E = "Present"; F = "Absent";  G = "Present"; I = "Absent";  J = "Present" }
    3 = @{ E = "Present"; F = "Present"; G = "Present"; H = "Absent";  I = "Present"; J = "Absent" }
    4 = @{ F = "Present"; H = "Absent";  I = "Present"; J = "Present" }
    5 = @{ E = "Present"; F = "Absent";  G = "Absent";  H = "Present"; J = "Absent" }
    6 = @{ E = "Present"; F = "Present"; I = "Absent";  J = "Present" }
}

# Only columns in the Day1..Day40 attendance block are considered.
$columns = @("E", "F", "G", "H", "I", "J", "K", "L", "M", "N")

# Columns that must always be (re)written for a given row, even if they
# already hold a value (e.g. Day1 gets recomputed / corrected).
$alwaysWrite = @{
    2 = @("E")
    3 = @("E")
    5 = @("E")
    6 = @("E")
}

foreach ($row in $attendance.Keys) {
    $rowValues = $attendance[$row]
    $forced = $alwaysWrite[$row]
    foreach ($col in $columns) {
        if (-not $rowValues.ContainsKey($col)) {
            continue
        }
        $cell = $ws.Range("$col$row")
        $isForced = ($forced -ne $null) -and ($forced -contains $col)
        # Automatically find empty columns/cells and fill those in; also
        # overwrite columns flagged for a forced refresh.
        if ($isForced -or $cell.Text -eq "") {
            $cell.Value = $rowValues[$col]
        }
    }
}
